$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = 'Última actualización: 20:11:56'
$ws1.Cells.Item(3,1).Value = 'Total filas: 132'

$rows1 = @{
  87 = @('18:10:41', '19:31', '27_EL RETIRO', 81, 'LP1912')
  88 = @('17:35:09', '19:31', '215_EL PELIGRO', 116, 'LP1912')
  89 = @('18:44:34', '19:33', '23_HERNANDEZ', 49, 'LP1912')
  90 = @('18:31:18', '19:34', '23_HERNANDEZ', 63, 'LP1912')
  91 = @('19:11:59', '19:38', '16_SANTA ANA', 27, 'LP1912')
  92 = @('19:11:59', '19:39', '17X38_ROMERO', 28, 'LP1912')
  93 = @('19:11:59', '19:40', '16_SANTA ANA', 29, 'LP1912')
  94 = @('17:47:22', '19:40', '17X38_ROMERO', 113, 'LP1912')
  95 = @('19:35:31', '19:41', '16_SANTA ANA', 6, 'LP1912')
  96 = @('19:11:59', '19:43', '11_ETCHEVERRY', 32, 'LP1912')
  97 = @('17:47:22', '19:44', '11_ETCHEVERRY', 117, 'LP1912')
  98 = @('18:31:18', '19:46', '11_ETCHEVERRY', 75, 'LP1912')
  99 = @('19:47:58', '19:47', '16_SANTA ANA', 0, 'LP1912')
  100 = @('19:11:59', '19:50', '81_EL PELIGRO', 39, 'LP1912')
  101 = @('17:54:43', '19:51', '81_EL PELIGRO', 117, 'LP1912')
  102 = @('19:54:49', '19:54', '16_SANTA ANA', 0, 'LP1912')
  103 = @('18:10:41', '19:58', '14X44_ABASTO', 108, 'LP1912')
  104 = @('18:31:18', '19:59', '14X44_ABASTO', 88, 'LP1912')
  105 = @('18:10:41', '20:00', '215C_EL PATO', 110, 'LP1912')
  106 = @('19:47:58', '20:00', '16_SANTA ANA', 13, 'LP1912')
  107 = @('18:31:18', '20:01', '215C_EL PATO', 90, 'LP1912')
  108 = @('19:47:58', '20:02', '17X38_ROMERO', 15, 'LP1912')
  109 = @('19:11:59', '20:04', '23_HERNANDEZ', 53, 'LP1912')
  110 = @('19:47:58', '20:09', '23_HERNANDEZ', 22, 'LP1912')
  111 = @('19:35:31', '20:10', '23_HERNANDEZ', 35, 'LP1912')
  112 = @('20:11:56', '20:12', '16_SANTA ANA', 1, 'LP1912')
  113 = @('20:11:56', '20:12', '11_ETCHEVERRY', 1, 'LP1912')
  114 = @('18:31:18', '20:13', '11_ETCHEVERRY', 62, 'LP1912')
  115 = @('18:31:18', '20:14', '11_ETCHEVERRY', 103, 'LP1912')
  116 = @('19:11:59', '20:25', '15_ABASTO', 74, 'LP1912')
  117 = @('18:31:18', '20:26', '15_ABASTO', 115, 'LP1912')
  118 = @('18:44:34', '20:28', '10_OLMOS', 104, 'LP1912')
  119 = @('18:31:18', '20:29', '10_OLMOS', 118, 'LP1912')
  120 = @('20:11:56', '20:35', '16_SANTA ANA', 24, 'LP1912')
  121 = @('19:11:59', '20:43', '215B_EL PATO', 92, 'LP1912')
  122 = @('19:11:59', '20:44', '17X38_ROMERO', 93, 'LP1912')
  123 = @('18:52:04', '20:44', '215B_EL PATO', 112, 'LP1912')
  124 = @('18:52:04', '20:45', '17X38_ROMERO', 113, 'LP1912')
  125 = @('20:11:56', '20:49', '23_HERNANDEZ', 38, 'LP1912')
  126 = @('19:54:49', '20:50', '23_HERNANDEZ', 56, 'LP1912')
  127 = @('19:35:31', '20:52', '23_HERNANDEZ', 77, 'LP1912')
  128 = @('19:54:49', '20:56', '27_EL RETIRO', 62, 'LP1912')
  129 = @('19:11:59', '21:01', '215A_EL PATO', 110, 'LP1912')
  130 = @('19:11:59', '21:02', '27_EL RETIRO', 111, 'LP1912')
  131 = @('19:47:58', '21:06', '27_EL RETIRO', 79, 'LP1912')
  132 = @('19:35:31', '21:10', '27_EL RETIRO', 95, 'LP1912')
  133 = @('19:35:31', '21:23', '10_OLMOS', 108, 'LP1912')
  134 = @('20:11:56', '21:34', '23_HERNANDEZ', 83, 'LP1912')
  135 = @('20:11:56', '21:48', '11_ETCHEVERRY', 97, 'LP1912')
  136 = @('19:54:49', '21:49', '11_ETCHEVERRY', 115, 'LP1912')
  137 = @('20:11:56', '21:55', '84_COLONIA URQUIZA-ESC 49', 104, 'LP1912')
}
foreach ($r in $rows1.Keys) {
  $v = $rows1[$r]
  $ws1.Cells.Item([int]$r, 1).Value = $v[0]
  $ws1.Cells.Item([int]$r, 2).Value = $v[1]
  $ws1.Cells.Item([int]$r, 3).Value = $v[2]
  $ws1.Cells.Item([int]$r, 4).Value = $v[3]
  $ws1.Cells.Item([int]$r, 5).Value = $v[4]
}

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = 'Última actualización: 20:11:56'

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = 'Última actualización: 20:11:56'
$ws3.Cells.Item(3,1).Value = 'Total filas: 15'

$rows3 = @{
  18 = @('20:11:56', '20:43', '215C_LA PLATA', 32, 'L6203')
  19 = @('19:47:58', '21:33', '215C_LA PLATA', 106, 'L6203')
  20 = @('19:35:31', '21:34', '215C_LA PLATA', 119, 'L6203')
}
foreach ($r in $rows3.Keys) {
  $v = $rows3[$r]
  $ws3.Cells.Item([int]$r, 1).Value = $v[0]
  $ws3.Cells.Item([int]$r, 2).Value = $v[1]
  $ws3.Cells.Item([int]$r, 3).Value = $v[2]
  $ws3.Cells.Item([int]$r, 4).Value = $v[3]
  $ws3.Cells.Item([int]$r, 5).Value = $v[4]
}

Write-Output "edit complete"
